$d = $word.ActiveDocument

# --- Edit 1: "Burrower" trait -> "Burrowing" trait, with its description
#     split across four separate runs and the distance numbers updated. ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("Burrower. The mole rat has a burrowing speed of 15 ft. It costs the mole rat 5 feet of movement to begin burrowing within a pre-existing entry point.")
if (-not $found1) {
    throw "Could not find the Burrower paragraph text"
}
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="62AB7803" w14:textId="4399C8D3" w:rsidR="005B3423" w:rsidRDefault="005B3423" w:rsidP="00C57E9C"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Burrowing. </w:t></w:r><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:t>mole rat</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>has a burrowing speed of 10 feet through loose earth and 0 feet through solid rock and metal.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Edit 2: insert a <w:lastRenderedPageBreak/> before the "Rabies is
#     rampant..." run, keeping the rest of the paragraph's runs intact. ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("Rabies is rampant among mole rats, creating wild and unpredictable members of the species. Such is their aggression that there are tall tales of rabid mole rats even charging down confused deathclaws. Some enterprising souls have realized this makes rabid mole rats the perfect delivery system for land mines. Raiders and other groups capture rabid mole rats, strap mines (or even remote-detonated explosives) on their backs and release them from cages to explode spectacularly.")
if (-not $found2) {
    throw "Could not find the Rabies paragraph text"
}
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5C860CE1" w14:textId="054D360D" w:rsidR="00F85F36" w:rsidRDefault="009434CC" w:rsidP="00C527CE"><w:r><w:lastRenderedPageBreak/><w:t>Rabies is rampant among mole rats, creating wild and unpredictable members of the species. Such is their aggression that there are tall tales of rabid mole rats</w:t></w:r><w:r w:rsidR="00DB31D3"><w:t xml:space="preserve"> even</w:t></w:r><w:r><w:t xml:space="preserve"> charging down confused deathclaws.</w:t></w:r><w:r w:rsidR="00F85F36"><w:t xml:space="preserve"> Some enterprising souls have realized this makes rabid mole rats the perfect delivery system for land mines. Raiders and other groups capture rabid mole rats, strap mines (or even remote-detonated explosives) </w:t></w:r><w:r w:rsidR="00186D0B"><w:t xml:space="preserve">on their backs </w:t></w:r><w:r w:rsidR="00F85F36"><w:t xml:space="preserve">and release them from cages </w:t></w:r><w:r w:rsidR="001C223A"><w:t>to explode spectacularly.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Output "Applied both edits"
